# [java->qt] - update class doc
# Reproduces (via the Excel object model) the hand-authored change:
#  - header font bumped 11 -> 15 (still bold), header row a bit taller
#  - the package/class table is repainted as a green/red status board:
#      * Package / class columns -> green
#      * %completed -> green only when 100, otherwise red
#      * ready for testing -> green when "yes", otherwise red
#      * tested -> green when "yes", otherwise red
#  - util/util.h's %completed corrected from 100 to 80
#  - cosmetic: column widths re-fitted, selection moved, row 2 heightened

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$GREEN = 5296274   # RGB(146, 208, 80)  -> FF92D050
$RED   = 255       # RGB(255,   0,  0)  -> FFFF0000

# ---- data correction -------------------------------------------------
$ws.Range("D9").Value = 80

# ---- header row (B2:F2) ---------------------------------------------
$ws.Range("B2:F2").Font.Bold = $true
$ws.Range("B2:F2").Font.Size = 15
$ws.Rows(2).RowHeight = 19.5

# ---- status colouring for the data rows (3-9) ------------------------
# Package / class columns are always green.
$ws.Range("B3:C9").Interior.Color = $GREEN

for ($row = 3; $row -le 9; $row++) {

    $completed = $ws.Cells.Item($row, 4).Value2
    if ($completed -eq 100) {
        $ws.Cells.Item($row, 4).Interior.Color = $GREEN
    } else {
        $ws.Cells.Item($row, 4).Interior.Color = $RED
    }

    $ready = $ws.Cells.Item($row, 5).Value2
    if ($ready -eq "yes") {
        $ws.Cells.Item($row, 5).Interior.Color = $GREEN
    } else {
        $ws.Cells.Item($row, 5).Interior.Color = $RED
    }

    $tested = $ws.Cells.Item($row, 6).Value2
    if ($tested -eq "yes") {
        $ws.Cells.Item($row, 6).Interior.Color = $GREEN
    } else {
        $ws.Cells.Item($row, 6).Interior.Color = $RED
    }
}

# ---- cosmetics: columns / selection -----------------------------------
$ws.Columns("A:A").ColumnWidth = 2.142857142857143
$ws.Columns("B:F").AutoFit()

$ws.Range("H6").Select()
